$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.730.84"
$ws.Range("D3").Value = "2.042.97"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.51"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.08"
$ws.Range("E7").Value = "  +1.47%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -1.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0832"
$ws.Range("E10").Value = "  +2.40%  "
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("D12").Value = "2.345.52"
$ws.Range("E12").Value = "  +0.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.36"
$ws.Range("E13").Value = "  -1.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.44"
$ws.Range("E14").Value = "  +2.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.50"
$ws.Range("E15").Value = "  +6.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.763"
$ws.Range("E16").Value = "  +0.99%  "
$ws.Range("D17").Value = "2.042.84"
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("D18").Value = "37.692.13"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.31"
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.90"
$ws.Range("E20").Value = "  -2.09%  "
$ws.Range("E21").Value = "  +0.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "222.03"
$ws.Range("E22").Value = "  -1.13%  "
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("E25").Value = "  +2.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.91"
$ws.Range("E26").Value = "  +2.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.25"
$ws.Range("E27").Value = "  +0.39%  "
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.74"
$ws.Range("E29").Value = "  -0.90%  "
$ws.Range("E30").Value = "  -0.56%  "
$ws.Range("E31").Value = "  -1.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.25"
$ws.Range("E32").Value = "  +8.08%  "
$ws.Range("E33").Value = "  -1.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.50"
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0602"
$ws.Range("E35").Value = "  +0.35%  "
$ws.Range("E36").Value = "  +2.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.34"
$ws.Range("E37").Value = "  +4.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.46"
$ws.Range("E38").Value = "  +7.18%  "
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.10"
$ws.Range("E40").Value = "  +7.51%  "
$ws.Range("D41").Value = "1.531.16"
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.56"
$ws.Range("E42").Value = "  +1.02%  "
$ws.Range("E43").Value = "  -0.99%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("E45").Value = "  -3.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0888"
$ws.Range("E46").Value = "  -3.02%  "
$ws.Range("E47").Value = "  +0.19%  "
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.94"
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.03"
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("D51").Value = "2.234.24"
$ws.Range("E51").Value = "  +0.69%  "
